$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure B:E columns are formatted as text so values like "68.064.76" are not
# reinterpreted as numbers/dates by Excel, matching the inlineStr cell type in the source.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.064.76"
$ws.Range("E2").Value = "  +7.60%  "
$ws.Range("D3").Value = "3.626.77"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "418.46"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "130.23"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "0.658"
$ws.Range("E7").Value = "  +5.13%  "
$ws.Range("D8").Value = "3.613.25"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.762"
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.197"
$ws.Range("E11").Value = "  +28.34%  "
$ws.Range("D12").Value = "0.0000430"
$ws.Range("E12").Value = "  +96.15%  "
$ws.Range("D13").Value = "42.11"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "9.86"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "4.205.38"
$ws.Range("E15").Value = "  +4.55%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "20.06"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.587.73"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "1.12"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "67.958.44"
$ws.Range("E20").Value = "  +7.40%  "
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "460.28"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "90.51"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "3.07"
$ws.Range("E24").Value = "  -5.77%  "
$ws.Range("D25").Value = "13.28"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "10.07"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").Value = "3.30"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "35.51"
$ws.Range("E28").Value = "  +6.56%  "
$ws.Range("D29").Value = "4.96"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("D30").Value = "2.80"
$ws.Range("E30").Value = "  +5.49%  "
$ws.Range("D31").Value = "12.27"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("E32").Value = "  +5.32%  "
$ws.Range("D33").Value = "7.22"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "40.44"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.157"
$ws.Range("E35").Value = "  -5.83%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0812"
$ws.Range("E37").Value = "  +24.49%  "
$ws.Range("D38").Value = "55.99"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "0.0481"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "0.147"
$ws.Range("E40").Value = "  +8.09%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "147.83"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.60"
$ws.Range("E46").Value = "  +9.96%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "4.25"
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.168"
$ws.Range("E48").Value = "  +19.37%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "0.301"
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("D50").Value = "1.94"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "116.56"
$ws.Range("E51").Value = "  +10.26%  "
